$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.024354954476961
$ws.Cells.Item(2, 4).Value = 1.034299980542307
$ws.Cells.Item(2, 5).Value = 1.04515910461195
$ws.Cells.Item(2, 6).Value = 1.048352702062955
$ws.Cells.Item(2, 9).Value = 1.031118593597669
$ws.Cells.Item(2, 10).Value = 1.029530034393715
$ws.Cells.Item(2, 11).Value = 1.037099851446981
$ws.Cells.Item(2, 12).Value = 1.04792812813185
$ws.Cells.Item(2, 13).Value = 1.051112785467323
$ws.Cells.Item(2, 14).Value = 1.013886138705523

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.025173632948521
$ws.Cells.Item(3, 4).Value = 1.034920959601419
$ws.Cells.Item(3, 5).Value = 1.046030302594381
$ws.Cells.Item(3, 6).Value = 1.049189424180054
$ws.Cells.Item(3, 9).Value = 1.031212782043853
$ws.Cells.Item(3, 10).Value = 1.029988260586713
$ws.Cells.Item(3, 11).Value = 1.037530491278758
$ws.Cells.Item(3, 12).Value = 1.048610554271379
$ws.Cells.Item(3, 13).Value = 1.05176147184657
$ws.Cells.Item(3, 14).Value = 1.014039099334389

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.025704014984549
$ws.Cells.Item(4, 4).Value = 1.035323173521545
$ws.Cells.Item(4, 5).Value = 1.046595396636453
$ws.Cells.Item(4, 6).Value = 1.049731912699386
$ws.Cells.Item(4, 9).Value = 1.031272492485597
$ws.Cells.Item(4, 10).Value = 1.030284753492766
$ws.Cells.Item(4, 11).Value = 1.037808851050103
$ws.Cells.Item(4, 12).Value = 1.049052858689809
$ws.Cells.Item(4, 13).Value = 1.052181639708323
$ws.Cells.Item(4, 14).Value = 1.014138033652294

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.025927139693643
$ws.Cells.Item(5, 4).Value = 1.035492357895878
$ws.Cells.Item(5, 5).Value = 1.046833288365657
$ws.Cells.Item(5, 6).Value = 1.049960229832595
$ws.Cells.Item(5, 9).Value = 1.031297298232057
$ws.Cells.Item(5, 10).Value = 1.030409395430429
$ws.Cells.Item(5, 11).Value = 1.03792580197897
$ws.Cells.Item(5, 12).Value = 1.049238976359636
$ws.Cells.Item(5, 13).Value = 1.052358378191698
$ws.Cells.Item(5, 14).Value = 1.014179615256172

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.025964612184827
$ws.Cells.Item(6, 4).Value = 1.035520770133695
$ws.Cells.Item(6, 5).Value = 1.046873250505178
$ws.Cells.Item(6, 6).Value = 1.049998580176574
$ws.Cells.Item(6, 9).Value = 1.031301445818168
$ws.Cells.Item(6, 10).Value = 1.030430323107585
$ws.Cells.Item(6, 11).Value = 1.037945434320626
$ws.Cells.Item(6, 12).Value = 1.049270236425901
$ws.Cells.Item(6, 13).Value = 1.052388059170704
$ws.Cells.Item(6, 14).Value = 1.014186596368068

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.025706995792271
$ws.Cells.Item(7, 4).Value = 1.035325433805836
$ws.Cells.Item(7, 5).Value = 1.046598574078114
$ws.Cells.Item(7, 6).Value = 1.049734962484797
$ws.Cells.Item(7, 9).Value = 1.031272825107171
$ws.Cells.Item(7, 10).Value = 1.030286418980488
$ws.Cells.Item(7, 11).Value = 1.037810414036888
$ws.Cells.Item(7, 12).Value = 1.049055344925945
$ws.Cells.Item(7, 13).Value = 1.052184000906113
$ws.Cells.Item(7, 14).Value = 1.014138589309103

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.024631496696591
$ws.Cells.Item(8, 4).Value = 1.034509759597418
$ws.Cells.Item(8, 5).Value = 1.045453245063989
$ws.Cells.Item(8, 6).Value = 1.048635252768388
$ws.Cells.Item(8, 9).Value = 1.031150680494862
$ws.Cells.Item(8, 10).Value = 1.029684895407124
$ws.Cells.Item(8, 11).Value = 1.03724544803454
$ws.Cells.Item(8, 12).Value = 1.048158605289455
$ws.Cells.Item(8, 13).Value = 1.051331923177349
$ws.Cells.Item(8, 14).Value = 1.013937840767668

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.022741323306148
$ws.Cells.Item(9, 4).Value = 1.033075575313782
$ws.Cells.Item(9, 5).Value = 1.043445613288895
$ws.Cells.Item(9, 6).Value = 1.046705733086719
$ws.Cells.Item(9, 9).Value = 1.030926011328064
$ws.Cells.Item(9, 10).Value = 1.02862491279935
$ws.Cells.Item(9, 11).Value = 1.036247722545937
$ws.Cells.Item(9, 12).Value = 1.046584092297971
$ws.Cells.Item(9, 13).Value = 1.049833778380763
$ws.Cells.Item(9, 14).Value = 1.0135837992977

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.021484664675099
$ws.Cells.Item(10, 4).Value = 1.032121673549663
$ws.Cells.Item(10, 5).Value = 1.042114425088594
$ws.Cells.Item(10, 6).Value = 1.045425093128876
$ws.Cells.Item(10, 9).Value = 1.030769928774608
$ws.Cells.Item(10, 10).Value = 1.027918318985536
$ws.Cells.Item(10, 11).Value = 1.035581186398157
$ws.Cells.Item(10, 12).Value = 1.045538318851237
$ws.Cells.Item(10, 13).Value = 1.048837351521994
$ws.Cells.Item(10, 14).Value = 1.013347599648114

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.020941358537198
$ws.Cells.Item(11, 4).Value = 1.031709173926484
$ws.Cells.Item(11, 5).Value = 1.041539746157683
$ws.Cells.Item(11, 6).Value = 1.044871940721017
$ws.Cells.Item(11, 9).Value = 1.030700856784136
$ws.Cells.Item(11, 10).Value = 1.027612386529076
$ws.Cells.Item(11, 11).Value = 1.035292257516906
$ws.Cells.Item(11, 12).Value = 1.045086433715589
$ws.Cells.Item(11, 13).Value = 1.048406463140642
$ws.Cells.Item(11, 14).Value = 1.013245287622943

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.020739677633024
$ws.Cells.Item(12, 4).Value = 1.031556037161253
$ws.Cells.Item(12, 5).Value = 1.041326547388613
$ws.Cells.Item(12, 6).Value = 1.044666683643327
$ws.Cells.Item(12, 9).Value = 1.030674977602737
$ws.Cells.Item(12, 10).Value = 1.027498754981441
$ws.Cells.Item(12, 11).Value = 1.035184890635362
$ws.Cells.Item(12, 12).Value = 1.044918726619992
$ws.Cells.Item(12, 13).Value = 1.048246499523178
$ws.Cells.Item(12, 14).Value = 1.013207279439251

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.020782933089061
$ws.Cells.Item(13, 4).Value = 1.031588881687159
$ws.Cells.Item(13, 5).Value = 1.041372267334056
$ws.Cells.Item(13, 6).Value = 1.044710702521709
$ws.Cells.Item(13, 9).Value = 1.03068053884501
$ws.Cells.Item(13, 10).Value = 1.027523129069718
$ws.Cells.Item(13, 11).Value = 1.035207923247472
$ws.Cells.Item(13, 12).Value = 1.044954693855701
$ws.Cells.Item(13, 13).Value = 1.048280808274342
$ws.Cells.Item(13, 14).Value = 1.013215432539788

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.020924684926428
$ws.Cells.Item(14, 4).Value = 1.031696513869876
$ws.Cells.Item(14, 5).Value = 1.041522117712893
$ws.Cells.Item(14, 6).Value = 1.044854969844663
$ws.Cells.Item(14, 9).Value = 1.030698722143474
$ws.Cells.Item(14, 10).Value = 1.027602993592565
$ws.Cells.Item(14, 11).Value = 1.035283383463854
$ws.Cells.Item(14, 12).Value = 1.04507256806448
$ws.Cells.Item(14, 13).Value = 1.048393238702928
$ws.Cells.Item(14, 14).Value = 1.013242145952729

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.021012039815695
$ws.Cells.Item(15, 4).Value = 1.031762840792746
$ws.Cells.Item(15, 5).Value = 1.041614480330575
$ws.Cells.Item(15, 6).Value = 1.044943885374407
$ws.Cells.Item(15, 9).Value = 1.030709895974776
$ws.Cells.Item(15, 10).Value = 1.027652201513318
$ws.Cells.Item(15, 11).Value = 1.035329870962368
$ws.Cells.Item(15, 12).Value = 1.045145213275991
$ws.Cells.Item(15, 13).Value = 1.048462522448459
$ws.Cells.Item(15, 14).Value = 1.013258604327245

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.021520739857637
$ws.Cells.Item(16, 4).Value = 1.032149061460737
$ws.Cells.Item(16, 5).Value = 1.042152601348597
$ws.Cells.Item(16, 6).Value = 1.045461833138531
$ws.Cells.Item(16, 9).Value = 1.030774481578273
$ws.Cells.Item(16, 10).Value = 1.027938623400367
$ws.Cells.Item(16, 11).Value = 1.035600355152881
$ws.Cells.Item(16, 12).Value = 1.045568328963468
$ws.Cells.Item(16, 13).Value = 1.048865960342645
$ws.Cells.Item(16, 14).Value = 1.013354389041471

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.021840058787985
$ws.Cells.Item(17, 4).Value = 1.032391475047773
$ws.Cells.Item(17, 5).Value = 1.042490616196693
$ws.Cells.Item(17, 6).Value = 1.045787097168474
$ws.Cells.Item(17, 9).Value = 1.030814596723084
$ws.Cells.Item(17, 10).Value = 1.028118296504954
$ws.Cells.Item(17, 11).Value = 1.035769939568712
$ws.Cells.Item(17, 12).Value = 1.045833991559101
$ws.Cells.Item(17, 13).Value = 1.049119180370513
$ws.Cells.Item(17, 14).Value = 1.013414462950416

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.022026392602818
$ws.Cells.Item(18, 4).Value = 1.032532923357968
$ws.Cells.Item(18, 5).Value = 1.042687941743224
$ws.Cells.Item(18, 6).Value = 1.045976950524786
$ws.Cells.Item(18, 9).Value = 1.030837851668821
$ws.Cells.Item(18, 10).Value = 1.028223099337761
$ws.Cells.Item(18, 11).Value = 1.035868824867513
$ws.Cells.Item(18, 12).Value = 1.045989038800472
$ws.Cells.Item(18, 13).Value = 1.049266934308616
$ws.Cells.Item(18, 14).Value = 1.013449499538458

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.022089941218832
$ws.Cells.Item(19, 4).Value = 1.032581162434356
$ws.Cells.Item(19, 5).Value = 1.042755252954934
$ws.Cells.Item(19, 6).Value = 1.046041707999451
$ws.Cells.Item(19, 9).Value = 1.030845756637233
$ws.Cells.Item(19, 10).Value = 1.028258834789022
$ws.Cells.Item(19, 11).Value = 1.035902536996484
$ws.Cells.Item(19, 12).Value = 1.046041921243379
$ws.Cells.Item(19, 13).Value = 1.049317323856353
$ws.Cells.Item(19, 14).Value = 1.013461445500477

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.021805790551638
$ws.Cells.Item(20, 4).Value = 1.032365460913109
$ws.Cells.Item(20, 5).Value = 1.042454333093217
$ws.Cells.Item(20, 6).Value = 1.045752185711845
$ws.Cells.Item(20, 9).Value = 1.030810307589838
$ws.Cells.Item(20, 10).Value = 1.028099019008151
$ws.Cells.Item(20, 11).Value = 1.035751747902435
$ws.Cells.Item(20, 12).Value = 1.045805479078534
$ws.Cells.Item(20, 13).Value = 1.049092006558705
$ws.Cells.Item(20, 14).Value = 1.013408017947845

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.020882939018886
$ws.Cells.Item(21, 4).Value = 1.031664816541655
$ws.Cells.Item(21, 5).Value = 1.041477983246809
$ws.Cells.Item(21, 6).Value = 1.044812480945673
$ws.Cells.Item(21, 9).Value = 1.030693373759419
$ws.Cells.Item(21, 10).Value = 1.027579475322114
$ws.Cells.Item(21, 11).Value = 1.035261163564881
$ws.Cells.Item(21, 12).Value = 1.04503785308234
$ws.Cells.Item(21, 13).Value = 1.048360128314351
$ws.Cells.Item(21, 14).Value = 1.013234279651092

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.020303442328607
$ws.Cells.Item(22, 4).Value = 1.031224780334763
$ws.Cells.Item(22, 5).Value = 1.040865633510362
$ws.Cells.Item(22, 6).Value = 1.044222857761813
$ws.Cells.Item(22, 9).Value = 1.030618564279501
$ws.Cells.Item(22, 10).Value = 1.027252849395963
$ws.Cells.Item(22, 11).Value = 1.034952449248931
$ws.Cells.Item(22, 12).Value = 1.044556045290426
$ws.Cells.Item(22, 13).Value = 1.047900474587982
$ws.Cells.Item(22, 14).Value = 1.013125015126694

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.020610574094722
$ws.Cells.Item(23, 4).Value = 1.031458005081429
$ws.Cells.Item(23, 5).Value = 1.041190106927226
$ws.Cells.Item(23, 6).Value = 1.044535313118903
$ws.Cells.Item(23, 9).Value = 1.030658344115517
$ws.Cells.Item(23, 10).Value = 1.027425996535982
$ws.Cells.Item(23, 11).Value = 1.035116129152483
$ws.Cells.Item(23, 12).Value = 1.044811381560165
$ws.Cells.Item(23, 13).Value = 1.048144097075809
$ws.Cells.Item(23, 14).Value = 1.013182940850235

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.021821274644102
$ws.Cells.Item(24, 4).Value = 1.032377215418624
$ws.Cells.Item(24, 5).Value = 1.042470727349514
$ws.Cells.Item(24, 6).Value = 1.045767960287671
$ws.Cells.Item(24, 9).Value = 1.030812246108123
$ws.Cells.Item(24, 10).Value = 1.028107729671444
$ws.Cells.Item(24, 11).Value = 1.035759968027546
$ws.Cells.Item(24, 12).Value = 1.045818362361028
$ws.Cells.Item(24, 13).Value = 1.049104285064917
$ws.Cells.Item(24, 14).Value = 1.013410930177979

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.023229376428582
$ws.Cells.Item(25, 4).Value = 1.033445962987536
$ws.Cells.Item(25, 5).Value = 1.043963367910334
$ws.Cells.Item(25, 6).Value = 1.04720356313398
$ws.Cells.Item(25, 9).Value = 1.030985207611323
$ws.Cells.Item(25, 10).Value = 1.028898938442251
$ws.Cells.Item(25, 11).Value = 1.036505908764478
$ws.Cells.Item(25, 12).Value = 1.046990461131581
$ws.Cells.Item(25, 13).Value = 1.050220680303958
$ws.Cells.Item(25, 14).Value = 1.013675359840425

